# Scheduled market-data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H, I, J, K, L, M, N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# with freshly pulled marketboard values.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3927.862
$ws.Range("I62").Value = 4753.4116
$ws.Range("K62").Value = 4753.4116
$ws.Range("M62").Value = -4129.4116
$ws.Range("H65").Value = 3927.862
$ws.Range("I65").Value = 4753.4116
$ws.Range("K65").Value = 23767.058
$ws.Range("M65").Value = -20647.058
$ws.Range("H92").Value = 417.2143
$ws.Range("J92").Value = 499.75
$ws.Range("L92").Value = 499.75
$ws.Range("N92").Value = -2995.75
$ws.Range("H132").Value = 23844.605
$ws.Range("I132").Value = 3586.7576
$ws.Range("J132").Value = 90695.5
$ws.Range("K132").Value = 10760.2728
$ws.Range("L132").Value = 272086.5
$ws.Range("M132").Value = -8230.272799999999
$ws.Range("N132").Value = -277146.5
$ws.Range("H137").Value = 1675470.5
$ws.Range("I137").Value = 2564968
$ws.Range("K137").Value = 7694904
$ws.Range("M137").Value = -7692354
$ws.Range("H138").Value = 1387.74
$ws.Range("I138").Value = 685.57574
$ws.Range("J138").Value = 1733.582
$ws.Range("K138").Value = 2056.72722
$ws.Range("L138").Value = 5200.746
$ws.Range("M138").Value = 3083.27278
$ws.Range("N138").Value = -15480.746

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12224.549
$ws.Range("I32").Value = 11661.813
$ws.Range("K32").Value = 11661.813
$ws.Range("M32").Value = -11374.813
$ws.Range("H61").Value = 2495.25
$ws.Range("I61").Value = 1832.0834
$ws.Range("K61").Value = 1832.0834
$ws.Range("M61").Value = -1620.0834
$ws.Range("H74").Value = 1376.8541
$ws.Range("I74").Value = 1242.8422
$ws.Range("J74").Value = 1886.1
$ws.Range("K74").Value = 1242.8422
$ws.Range("L74").Value = 1886.1
$ws.Range("M74").Value = -368.8422
$ws.Range("N74").Value = -3634.1
$ws.Range("H77").Value = 1376.8541
$ws.Range("I77").Value = 1242.8422
$ws.Range("J77").Value = 1886.1
$ws.Range("K77").Value = 6214.211
$ws.Range("L77").Value = 9430.5
$ws.Range("M77").Value = -1846.211
$ws.Range("N77").Value = -18166.5
$ws.Range("H88").Value = 17871048
$ws.Range("I88").Value = 28575064
$ws.Range("K88").Value = 28575064
$ws.Range("M88").Value = -28574658
$ws.Range("H91").Value = 17871048
$ws.Range("I91").Value = 28575064
$ws.Range("K91").Value = 28575064
$ws.Range("M91").Value = -28573660
$ws.Range("H102").Value = 40947.168
$ws.Range("I102").Value = 3783
$ws.Range("K102").Value = 3783
$ws.Range("M102").Value = -2161
$ws.Range("H132").Value = 15153533
$ws.Range("I132").Value = 21740426
$ws.Range("J132").Value = 3679.6
$ws.Range("K132").Value = 65221278
$ws.Range("L132").Value = 11038.8
$ws.Range("M132").Value = -65218748
$ws.Range("N132").Value = -16098.8
$ws.Range("H136").Value = 2495.25
$ws.Range("I136").Value = 1832.0834
$ws.Range("K136").Value = 5496.2502
$ws.Range("M136").Value = -2946.2502

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2087.4666
$ws.Range("I107").Value = 1800.5238
$ws.Range("K107").Value = 1800.5238
$ws.Range("M107").Value = 119.4762000000001
$ws.Range("H134").Value = 3248.9048
$ws.Range("I134").Value = 2551.4546
$ws.Range("J134").Value = 3623.1462
$ws.Range("K134").Value = 7654.3638
$ws.Range("L134").Value = 10869.4386
$ws.Range("M134").Value = -5119.3638
$ws.Range("N134").Value = -15939.4386

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2632.04
$ws.Range("I31").Value = 873.5
$ws.Range("J31").Value = 3187.3684
$ws.Range("K31").Value = 873.5
$ws.Range("L31").Value = 3187.3684
$ws.Range("M31").Value = -578.5
$ws.Range("N31").Value = -3777.3684
$ws.Range("H34").Value = 2632.04
$ws.Range("I34").Value = 873.5
$ws.Range("J34").Value = 3187.3684
$ws.Range("K34").Value = 873.5
$ws.Range("L34").Value = 3187.3684
$ws.Range("M34").Value = -671.5
$ws.Range("N34").Value = -3591.3684
$ws.Range("H58").Value = 1940.1482
$ws.Range("I58").Value = 1614.45
$ws.Range("J58").Value = 2870.7144
$ws.Range("K58").Value = 1614.45
$ws.Range("L58").Value = 2870.7144
$ws.Range("M58").Value = -1411.45
$ws.Range("N58").Value = -3276.7144
$ws.Range("H132").Value = 63656.566
$ws.Range("I132").Value = 1301.6923
$ws.Range("J132").Value = 144717.9
$ws.Range("K132").Value = 3905.0769
$ws.Range("L132").Value = 434153.7
$ws.Range("M132").Value = -1375.0769
$ws.Range("N132").Value = -439213.7
$ws.Range("H134").Value = 740464.1
$ws.Range("I134").Value = 522663.78
$ws.Range("J134").Value = 1275065
$ws.Range("K134").Value = 1567991.34
$ws.Range("L134").Value = 3825195
$ws.Range("M134").Value = -1565456.34
$ws.Range("N134").Value = -3830265
$ws.Range("H136").Value = 1940.1482
$ws.Range("I136").Value = 1614.45
$ws.Range("J136").Value = 2870.7144
$ws.Range("K136").Value = 4843.35
$ws.Range("L136").Value = 8612.143199999999
$ws.Range("M136").Value = -2293.35
$ws.Range("N136").Value = -13712.1432

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 11467
$ws.Range("I122").Value = 359.33334
$ws.Range("K122").Value = 3234.00006
$ws.Range("M122").Value = -784.0000600000003
$ws.Range("H131").Value = 939.55
$ws.Range("J131").Value = 960.5789
$ws.Range("L131").Value = 2881.7367
$ws.Range("N131").Value = -12961.7367
$ws.Range("H141").Value = 166670990
$ws.Range("I141").Value = 333336260
$ws.Range("J141").Value = 5733
$ws.Range("K141").Value = 1000008780
$ws.Range("L141").Value = 17199
$ws.Range("M141").Value = -1000003600
$ws.Range("N141").Value = -27559

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4148.5713
$ws.Range("J102").Value = 3205.6
$ws.Range("L102").Value = 3205.6
$ws.Range("N102").Value = -6449.6
$ws.Range("H126").Value = 12381.1
$ws.Range("I126").Value = 28427.75
$ws.Range("J126").Value = 1683.3334
$ws.Range("K126").Value = 85283.25
$ws.Range("L126").Value = 5050.0002
$ws.Range("M126").Value = -82813.25
$ws.Range("N126").Value = -9990.0002
$ws.Range("H132").Value = 27029834
$ws.Range("I132").Value = 43479916
$ws.Range("K132").Value = 130439748
$ws.Range("M132").Value = -130437218

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2383.5073
$ws.Range("I132").Value = 1776.8679
$ws.Range("J132").Value = 4393
$ws.Range("K132").Value = 5330.6037
$ws.Range("L132").Value = 13179
$ws.Range("M132").Value = -2800.6037
$ws.Range("N132").Value = -18239
$ws.Range("H136").Value = 1735.3658
$ws.Range("I136").Value = 1353.0322
$ws.Range("J136").Value = 2920.6
$ws.Range("K136").Value = 4059.0966
$ws.Range("L136").Value = 8761.799999999999
$ws.Range("M136").Value = -1509.0966
$ws.Range("N136").Value = -13861.8

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4741.4116
$ws.Range("I81").Value = 1018.36365
$ws.Range("J81").Value = 11567
$ws.Range("K81").Value = 2036.7273
$ws.Range("L81").Value = 23134
$ws.Range("M81").Value = -975.7273
$ws.Range("N81").Value = -25256
$ws.Range("H84").Value = 4741.4116
$ws.Range("I84").Value = 1018.36365
$ws.Range("J84").Value = 11567
$ws.Range("K84").Value = 10183.6365
$ws.Range("L84").Value = 115670
$ws.Range("M84").Value = -4879.636500000001
$ws.Range("N84").Value = -126278
$ws.Range("H132").Value = 1209406.8
$ws.Range("I132").Value = 1553858.8
$ws.Range("J132").Value = 3824.5
$ws.Range("K132").Value = 4661576.4
$ws.Range("L132").Value = 11473.5
$ws.Range("M132").Value = -4659046.4
$ws.Range("N132").Value = -16533.5
$ws.Range("H136").Value = 805671.3
$ws.Range("I136").Value = 1167415.1
$ws.Range("J136").Value = 1796.1111
$ws.Range("K136").Value = 3502245.3
$ws.Range("L136").Value = 5388.3333
$ws.Range("M136").Value = -3499695.3
$ws.Range("N136").Value = -10488.3333
